# "new features to MaxTable"
#
# The first 7 columns (A:G) of the header row are reorganized: the three
# columns "Module PN", "Revision" and "Wire Type" move to the front
# (A, B, C) and "Customer Module/Option Function" moves to the very end of
# that block (column G). The columns carry their original widths with them.
# Everything from column H onward is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- New header text for columns A-G (content follows the moved columns) ----
$ws.Range("A1").Value = "Module PN"
$ws.Range("B1").Value = "Revision"
$ws.Range("C1").Value = "Wire Type"
$ws.Range("D1").Value = "Customer Module/Option"
$ws.Range("E1").Value = "Wire Customer Name"
$ws.Range("F1").Value = "Wire Internal Name"
$ws.Range("G1").Value = "Customer Module/Option Function"

# ---- Column widths travel with the moved header text ----
# ColumnWidth is stored internally in whole-pixel steps (1/6 "character"
# units here), so we solve for the input that lands exactly back on the
# desired stored width instead of feeding the target straight in (which
# would drift by a fraction of a character due to rounding).
function Set-ExactColumnWidth($col, $targetWidth) {
    $px = [math]::Round($targetWidth * 6)
    $compensated = ($px - 5) / 6.0
    $ws.Columns($col).ColumnWidth = $compensated
}

Set-ExactColumnWidth "A" 15.5703125
Set-ExactColumnWidth "B" 13.140625
Set-ExactColumnWidth "C" 14.7109375
Set-ExactColumnWidth "D" 28.85546875
Set-ExactColumnWidth "E" 24.85546875
Set-ExactColumnWidth "F" 23.28515625
Set-ExactColumnWidth "G" 37.28515625

# ---- Selection left on the (now relocated) "Customer Module/Option
# Function" column, as if the user had just finished moving it and it's
# still selected as a whole column ----
$ws.Columns("G").Select() | Out-Null
